$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Fix-Text($s) {
    $s = $s -replace '\bde\b', 'De'
    $s = $s -replace '\bdel\b', 'Del'
    $s = $s -replace '\bla\b', 'La'
    $s = $s -replace '\bel\b', 'El'
    $s = $s -replace '\blos\b', 'Los'
    $s = $s -replace '\blas\b', 'Las'
    $s = $s -replace '\by\b', 'Y'
    return $s
}

# 1) Rename header row (A1:D1) to the new short column names.
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2) Title-case the connector words ("de", "del", "la", "el", "los", "las", "y")
#    inside the state (column A) and municipality (column B) text cells, and
#    fix the one-ULP float drift on the percentage column for rows where the
#    count (column C) equals 4.
$lastRow = 908

for ($r = 2; $r -le $lastRow; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    if ($aCell.Text -ne "") {
        $aCell.Value = Fix-Text $aCell.Text
    }

    $bCell = $ws.Cells.Item($r, 2)
    if ($bCell.Text -ne "") {
        $bCell.Value = Fix-Text $bCell.Text
    }

    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2 -eq 4) {
        $ws.Cells.Item($r, 4).Value = 0.0009250693802035152
    }
}

# 3) Drop the trailing metadata/footer rows (909-914) and shrink the sheet
#    back down to the real data range (dimension becomes A1:D908).
$ws.Rows("909:914").Delete()
